# update and regroup Maximum_Subarray
#
# - Remove the "Maximum Subarray" entry (row 4) from the DP_Matrix sheet,
#   shifting the following rows up.
# - Add "Maximum Subarray" into the DP_Sequence sheet at row 6 (column A),
#   keeping the existing formatting of B6.
# - Update sheet view / active-tab state to match: DP_Matrix becomes the
#   selected tab (instead of Divide&Conquer), with new selections on the
#   DP_Matrix and DP_Sequence sheets.

$wb = $excel.ActiveWorkbook

$wsDivideConquer = $wb.Worksheets.Item("Divide&Conquer")
$wsDPMatrix = $wb.Worksheets.Item("DP_Matrix")
$wsDPSequence = $wb.Worksheets.Item("DP_Sequence")

# Remove "Maximum Subarray" row from DP_Matrix (row 4), shifting rows 5-9 up.
$wsDPMatrix.Rows(4).Delete()

# Add "Maximum Subarray" to DP_Sequence, row 6, column A (B6 keeps its style).
$wsDPSequence.Cells.Item(6, 1).Value = "Maximum Subarray"

# Update selections / active sheet state.
$wsDivideConquer.Activate()
$wsDivideConquer.Range("O16").Select()

$wsDPSequence.Activate()
$wsDPSequence.Range("F47").Select()

$wsDPMatrix.Activate()
$wsDPMatrix.Range("F10").Select()
